# This script applies a weekly data refresh to the "Plátano" (banana) price
# sheet: two new price records (for a new reporting date) are inserted at
# row 269, pushing the existing rows 269-376 down to rows 271-378.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at row 269 (this shifts old rows 269:376 down to 271:378,
# carrying their formatting/styles along, same as native Excel row insertion).
$ws.Range("A269:A270").EntireRow.Insert()

# Column layout (same for every data row in this sheet):
#  A Mercado ID            (numeric)
#  B Mercado                (text)
#  C Región                 (text)
#  D Fecha                  (numeric / date serial)
#  E Codreg                 (numeric)
#  F Tipo                   (text)
#  G Producto ID             (numeric)
#  H Producto                (text)
#  I Categoría ID            (numeric)
#  J Categoría                (text)
#  K Variedad                 (text)
#  L Calidad                  (text)
#  M Volumen                  (numeric)
#  N Precio mínimo            (numeric)
#  O Precio máximo            (numeric)
#  P Precio promedio ponderado (numeric)
#  Q Unidad de comercialización (text)
#  R Origen                    (text)
#  S Precio $/Kg                (numeric)
#  T Kg / unidad                 (numeric)

$numericCols = @(1,4,5,7,9,13,14,15,16,19,20)

$row269 = @{
    1  = 1
    2  = "Agrícola del Norte S.A. de Arica"
    3  = "Arica y Parinacota"
    4  = 45009
    5  = 15
    6  = "Fruta"
    7  = 100108
    8  = "Tropicales y subtropicales"
    9  = 100108006
    10 = "Plátano"
    11 = "Sin especificar"
    12 = "Pintón"
    13 = 120
    14 = 24000
    15 = 25000
    16 = 24500
    17 = "`$/caja 20 kilos"
    18 = "Ecuador"
    19 = 1225
    20 = 20
}

$row270 = @{
    1  = 1
    2  = "Agrícola del Norte S.A. de Arica"
    3  = "Arica y Parinacota"
    4  = 45009
    5  = 15
    6  = "Fruta"
    7  = 100108
    8  = "Tropicales y subtropicales"
    9  = 100108006
    10 = "Plátano"
    11 = "Sin especificar"
    12 = "Verde"
    13 = 120
    14 = 22000
    15 = 23000
    16 = 22500
    17 = "`$/caja 20 kilos"
    18 = "Ecuador"
    19 = 1125
    20 = 20
}

foreach ($col in 1..20) {
    $ws.Cells.Item(269, $col).Value2 = $row269[$col]
    $ws.Cells.Item(270, $col).Value2 = $row270[$col]
}
